# "Generate Report for Handoff" - regenerate report with a freshly-generated
# handoff id (replacing the old 3cbda7d4... guid with 572d7470...) and new
# handoff timestamps; the not-yet-handed-back target/handback columns on the
# zh-cn / de-de sheets are reset to blank pending the next handback.

$wb = $excel.ActiveWorkbook

$oldId = "3cbda7d4-8939-4a3b-9ac8-dbec18455361"
$newId = "572d7470-b459-4fec-b471-e3302e8245ef"
$newHash = "f1f1b6447d7391c73d7303e2ec25a0f531cd51f3"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newId.md"
$ws1.Range("G2").Value = "2016-08-21 21:09:42"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6487e950b8f2032666a4ca5751b5de429d057ccd/e2e/$oldId.md"
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $url1, "", "", "e2e\$newId.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$urlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6487e950b8f2032666a4ca5751b5de429d057ccd/e2e/$oldId.md"

# Remove both hyperlinks (range-scoped delete removes every hyperlink on the
# sheet in this engine), then re-add only the one that survives (A2) - the
# "Latest Target File" hyperlink on I2 is gone now that it is blank again.
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlA2, "", "", "$newId.md")

$ws2.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-21 21:09:38"

# Latest Target File / Latest Handback File reset to blank (still typed as
# text, like the sheet's other blank cells) pending the next handback.
$ws2.Range("I2").Value = "'"
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = "'"
$ws2.Range("J2").Style = "Normal"

$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$urlA2de = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6487e950b8f2032666a4ca5751b5de429d057ccd/e2e/$oldId.md"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlA2de, "", "", "$newId.md")

$ws3.Range("G2").Value = "$newId.$newHash.de-de.xlf"

$ws3.Range("I2").Value = "'"
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = "'"
$ws3.Range("J2").Style = "Normal"

$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
